$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2191.875
$ws.Range("I106").Value = 1923.3334
$ws.Range("J106").Value = 2997.5
$ws.Range("K106").Value = 1923.3334
$ws.Range("L106").Value = 2997.5
$ws.Range("M106").Value = -1292.3334
$ws.Range("N106").Value = -4259.5
$ws.Range("H116").Value = 3876.25
$ws.Range("I116").Value = 4000
$ws.Range("K116").Value = 4000
$ws.Range("M116").Value = -558
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 405.4
$ws.Range("J5").Value = 325.66666
$ws.Range("L5").Value = 325.66666
$ws.Range("N5").Value = -549.66666
$ws.Range("H61").Value = 5567.357
$ws.Range("I61").Value = 3676.6365
$ws.Range("K61").Value = 3676.6365
$ws.Range("M61").Value = -3464.6365
$ws.Range("H110").Value = 62500690
$ws.Range("I110").Value = 62500690
$ws.Range("K110").Value = 62500690
$ws.Range("M110").Value = -62498645
$ws.Range("H132").Value = 5429.3076
$ws.Range("I132").Value = 3329.2
$ws.Range("J132").Value = 12429.667
$ws.Range("K132").Value = 9987.599999999999
$ws.Range("L132").Value = 37289.001
$ws.Range("M132").Value = -7457.599999999999
$ws.Range("N132").Value = -42349.001
$ws.Range("H136").Value = 5567.357
$ws.Range("I136").Value = 3676.6365
$ws.Range("K136").Value = 11029.9095
$ws.Range("M136").Value = -8479.9095

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 405.4
$ws.Range("J4").Value = 325.66666
$ws.Range("L4").Value = 325.66666
$ws.Range("N4").Value = -555.66666
$ws.Range("H105").Value = 500001250
$ws.Range("I105").Value = 1000000000
$ws.Range("K105").Value = 1000000000
$ws.Range("M105").Value = -999998253
$ws.Range("H107").Value = 2156.5557
$ws.Range("I107").Value = 2177.5881
$ws.Range("J107").Value = 1799
$ws.Range("K107").Value = 2177.5881
$ws.Range("L107").Value = 1799
$ws.Range("M107").Value = -257.5880999999999
$ws.Range("N107").Value = -5639
$ws.Range("H134").Value = 4221.6665
$ws.Range("I134").Value = 2499.348
$ws.Range("K134").Value = 7498.044
$ws.Range("M134").Value = -4963.044

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 373.92307
$ws.Range("I7").Value = 169.14285
$ws.Range("J7").Value = 612.8333
$ws.Range("K7").Value = 169.14285
$ws.Range("L7").Value = 612.8333
$ws.Range("M7").Value = -56.14285000000001
$ws.Range("N7").Value = -838.8333
$ws.Range("H58").Value = 4226.5
$ws.Range("I58").Value = 2727.5715
$ws.Range("K58").Value = 2727.5715
$ws.Range("M58").Value = -2524.5715
$ws.Range("H96").Value = 44000
$ws.Range("J96").Value = 44000
$ws.Range("L96").Value = 44000
$ws.Range("N96").Value = -49492
$ws.Range("H103").Value = 22470.875
$ws.Range("J103").Value = 79979
$ws.Range("L103").Value = 79979
$ws.Range("N103").Value = -82323
$ws.Range("H110").Value = 55000
$ws.Range("J110").Value = 55000
$ws.Range("L110").Value = 55000
$ws.Range("N110").Value = -63180
$ws.Range("H122").Value = 334166.66
$ws.Range("I122").Value = 500250.5
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 1500751.5
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -1498301.5
$ws.Range("N122").Value = -10897
$ws.Range("H136").Value = 4226.5
$ws.Range("I136").Value = 2727.5715
$ws.Range("K136").Value = 8182.7145
$ws.Range("M136").Value = -5632.7145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2165.4849
$ws.Range("J2").Value = 3559.15
$ws.Range("L2").Value = 21354.9
$ws.Range("N2").Value = -21580.9
$ws.Range("H5").Value = 735
$ws.Range("J5").Value = 802.5
$ws.Range("L5").Value = 2407.5
$ws.Range("N5").Value = -2631.5
$ws.Range("H88").Value = 4669.4443
$ws.Range("J88").Value = 7016
$ws.Range("L88").Value = 21048
$ws.Range("N88").Value = -21904
$ws.Range("H91").Value = 4669.4443
$ws.Range("J91").Value = 7016
$ws.Range("L91").Value = 21048
$ws.Range("N91").Value = -24012
$ws.Range("H95").Value = 7231.75
$ws.Range("J95").Value = 9610.799999999999
$ws.Range("L95").Value = 28832.4
$ws.Range("N95").Value = -32950.39999999999
$ws.Range("H135").Value = 735
$ws.Range("J135").Value = 802.5
$ws.Range("L135").Value = 7222.5
$ws.Range("N135").Value = -12292.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 375697.28
$ws.Range("J14").Value = 1693.75
$ws.Range("L14").Value = 1693.75
$ws.Range("N14").Value = -2029.75
$ws.Range("H19").Value = 8270.857
$ws.Range("I19").Value = 4500
$ws.Range("J19").Value = 11099
$ws.Range("K19").Value = 4500
$ws.Range("L19").Value = 11099
$ws.Range("M19").Value = -4212
$ws.Range("N19").Value = -11675
$ws.Range("H80").Value = 2107
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002
$ws.Range("H83").Value = 2107
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("M83").Value = -5008
$ws.Range("H97").Value = 823.4286
$ws.Range("I97").Value = 770.7273
$ws.Range("J97").Value = 1016.6667
$ws.Range("K97").Value = 770.7273
$ws.Range("L97").Value = 1016.6667
$ws.Range("M97").Value = -274.7273
$ws.Range("N97").Value = -2008.6667
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H102").Value = 1460.9
$ws.Range("I102").Value = 1156.875
$ws.Range("K102").Value = 1156.875
$ws.Range("M102").Value = 465.125
$ws.Range("H104").Value = 78000
$ws.Range("J104").Value = 78000
$ws.Range("L104").Value = 78000
$ws.Range("N104").Value = -84988
$ws.Range("H105").Value = 57666
$ws.Range("J105").Value = 57666
$ws.Range("L105").Value = 57666
$ws.Range("N105").Value = -64654
$ws.Range("H112").Value = 50293
$ws.Range("J112").Value = 50293
$ws.Range("L112").Value = 50293
$ws.Range("N112").Value = -52509
$ws.Range("H128").Value = 10000
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value = 94980
$ws.Range("J129").Value = 94980
$ws.Range("L129").Value = 94980
$ws.Range("N129").Value = -104980
$ws.Range("H132").Value = 6716.2593
$ws.Range("I132").Value = 4077.6875
$ws.Range("K132").Value = 12233.0625
$ws.Range("M132").Value = -9703.0625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7373.759
$ws.Range("J132").Value = 8065.4
$ws.Range("L132").Value = 24196.2
$ws.Range("N132").Value = -29256.2
$ws.Range("H136").Value = 4421.2
$ws.Range("I136").Value = 2497.9285
$ws.Range("J136").Value = 6869
$ws.Range("K136").Value = 7493.7855
$ws.Range("L136").Value = 20607
$ws.Range("M136").Value = -4943.7855
$ws.Range("N136").Value = -25707

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 28799.8
$ws.Range("J43").Value = 51999.5
$ws.Range("L43").Value = 51999.5
$ws.Range("N43").Value = -52297.5
$ws.Range("H113").Value = 529.1429000000001
$ws.Range("I113").Value = 475
$ws.Range("J113").Value = 664.5
$ws.Range("K113").Value = 1425
$ws.Range("L113").Value = 1993.5
$ws.Range("M113").Value = 745
$ws.Range("N113").Value = -6333.5
$ws.Range("H136").Value = 2518.6943
$ws.Range("I136").Value = 1904.7916
$ws.Range("J136").Value = 3746.5
$ws.Range("K136").Value = 5714.3748
$ws.Range("L136").Value = 11239.5
$ws.Range("M136").Value = -3164.3748
$ws.Range("N136").Value = -16339.5
